# Applies the "adjusted some texts to be more 'simple'. Changed explanation
# of humidity." commit to the document.
#
# Strategy: every edit below is applied via a precise Range (character-offset)
# replacement so we never depend on literal run boundaries. Paragraph counts
# are preserved by every single edit (no paragraph is added or removed
# anywhere in the document), so 1-based $d.Paragraphs(...) indices computed
# once up-front stay valid for the whole script.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$OldSub,
        [string]$NewSub
    )
    $doc = $word.ActiveDocument
    $p = $doc.Paragraphs($ParaIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($OldSub)
    if ($idx -lt 0) {
        throw "Replace-InParagraph: substring not found in paragraph $ParaIndex : [$OldSub] (actual text: [$full])"
    }
    $start = $p.Range.Start + $idx
    $end = $start + $OldSub.Length
    $r = $doc.Range($start, $end)
    $r.Text = $NewSub
}

function Replace-ParagraphRange {
    param(
        [int]$FirstParaIndex,
        [int]$LastParaIndex,
        [string]$NewText
    )
    $doc = $word.ActiveDocument
    $pStart = $doc.Paragraphs($FirstParaIndex).Range.Start
    $pEnd = $doc.Paragraphs($LastParaIndex).Range.End
    $r = $doc.Range($pStart, $pEnd)
    $r.Text = $NewText
}

# ---------------------------------------------------------------------------
# 1) "...are visited by sentient beings, and frogs at that!"
#    -> "...are visited by space frogs!"
# ---------------------------------------------------------------------------
Replace-InParagraph 8 "visited by sentient beings, and frogs at that!" "visited by space frogs!"

# ---------------------------------------------------------------------------
# 2) "...give these hapless frogs..." -> "...give these frogs..."
# ---------------------------------------------------------------------------
Replace-InParagraph 11 "give these hapless " "give these "

# ---------------------------------------------------------------------------
# 3) Humidity explanation rewrite (paragraphs 30-32 -> 3 new paragraphs)
# ---------------------------------------------------------------------------
$humidityText = "Next is the humidity readings of Earth.`r" + `
    "Humidity tells us how much water vapor is in the air. These water vapor comes from evaporation, and is dropped to new location as the air cools down.`r" + `
    "What you see on the map is the relative humidity in percentage. This is the amount of water in the air relative to the maximum amount of water vapor (moisture). "
Replace-ParagraphRange 30 32 $humidityText

# ---------------------------------------------------------------------------
# 4) Table cell: remove proofErr markers around "So" (text content unchanged)
# ---------------------------------------------------------------------------
Replace-InParagraph 59 "So the weather we will be experiencing here over a few days will be vastly different several months later." "So the weather we will be experiencing here over a few days will be vastly different several months later."

# ---------------------------------------------------------------------------
# 5) "...populate the place with as many frogs as possible. We do this by
#    deploying houses." -> "...populate the land with as many frogs as
#    possible. We do this by placing houses."
# ---------------------------------------------------------------------------
Replace-InParagraph 62 "populate the place with as many frogs as possible. We do this by deploying houses." "populate the land with as many frogs as possible. We do this by placing houses."

# ---------------------------------------------------------------------------
# 6) "You can deploy a house by pressing..." -> "You can place a house by
#    pressing..."
# ---------------------------------------------------------------------------
Replace-InParagraph 63 "You can deploy a house by pressing" "You can place a house by pressing"

# ---------------------------------------------------------------------------
# 7-10) Remove proofErr markers (spell-check artifacts); text is unchanged.
# ---------------------------------------------------------------------------
Replace-InParagraph 106 "(toggle to wind attr)" "(toggle to wind attr)"
Replace-InParagraph 119 "(toggle temp attr)" "(toggle temp attr)"
Replace-InParagraph 125 "(hide  gulf stream illustration)" "(hide  gulf stream illustration)"
Replace-InParagraph 147 "(hazzard)" "(hazzard)"
